# Update the DNA extraction example row (row 2) on the "2EXT04_DNA" sheet
# to reflect the new sample description / ontology references, per the
# commit "update DNA extraction example values and ontology ref for
# sample description".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2EXT04_DNA")

# Characteristic [bio entity]: "total RNA" -> "DNA"
$ws.Range("B2").Value = "DNA"

# Term Source REF (DPBO:0000012): "EFO" -> "BAO"
$ws.Range("C2").Value = "BAO"

# Term Accession Number (DPBO:0000012)
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/BAO_0000269"

# Parameter [biosource amount]: "200" -> "1"
# (force text so the numeric-looking value is stored like the original,
# not auto-converted to a number, then restore the cell's prior style)
$styleE2 = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1"
$ws.Range("E2").Style = $styleE2

# Unit: "milligram" -> "gram"
$ws.Range("F2").Value = "gram"

# Term Accession Number (DPBO:0000013): UO:0000022 -> UO:0000021
$ws.Range("H2").Value = "https://bioregistry.io/UO:0000021"

# Parameter [extraction method]
$ws.Range("I2").Value = "Macherey Nagel NucleoBond HMW DNA Kit"

# Component [extraction buffer]
$ws.Range("L2").Value = "Lysis buffer H1"

# Parameter [extraction buffer volume]: "200" -> "4"
$styleO2 = $ws.Range("O2").Style
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "4"
$ws.Range("O2").Style = $styleO2

# Unit : "microliter" -> "milliliter"
$ws.Range("P2").Value = "milliliter"

# Term Accession Number (DPBO:0000051): UO:0000101 -> UO:0000098
$ws.Range("R2").Value = "https://bioregistry.io/UO:0000098"
